$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing last row (row 6) down to row 7 by copying it,
# then overwrite row 6 with the new report entry.
$ws.Range("A6:H6").Copy($ws.Range("A7:H7"))

# Fill in the new row 6 with the newest weekly-report entry.
$ws.Range("A6").Value = "24/9/2012"
$ws.Range("B6").Value = "chỉnh sửa hoàn thiện SRS cá nhân"
$ws.Range("C6").Value = "1phaanf SRS"
$ws.Range("D6").Value = "Hoàn thành"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.5

# Match the taller row used for similar wrapped-text entries.
$ws.Rows.Item(6).RowHeight = 40.5

# Reflect where the user left the selection after editing.
$ws.Range("F7").Select()
